$d = $word.ActiveDocument

$replacements = @(
    @("2024-07-27 Saturday", "2024-07-28 Sunday"),
    @("290×8=2320", "872×6=5232"),
    @("147×7=1029", "708×6=4248"),
    @("993×9=8937", "875×2=1750"),
    @("789×7=5523", "451×7=3157"),
    @("596×9=5364", "344×5=1720"),
    @("221×3=663", "402×7=2814"),
    @("944×9=8496", "569×2=1138"),
    @("573×3=1719", "549×6=3294"),
    @("138×3=414", "146×2=292"),
    @("336×5=1680", "199×9=1791"),
    @("563×9=5067", "125×4=500"),
    @("990×2=1980", "744×8=5952"),
    @("867×9=7803", "692×2=1384"),
    @("289×8=2312", "801×4=3204"),
    @("812×3=2436", "679×7=4753"),
    @("520×2=1040", "496×2=992"),
    @("678×4=2712", "619×2=1238"),
    @("479×9=4311", "808×5=4040"),
    @("709×9=6381", "516×9=4644"),
    @("112×9=1008", "665×9=5985"),
    @("162×4=648", "251×8=2008"),
    @("203×9=1827", "633×9=5697"),
    @("168×8=1344", "588×8=4704"),
    @("699×4=2796", "825×9=7425"),
    @("538×4=2152", "880×7=6160")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $found = $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
    if (-not $found) {
        Write-Host "NOT FOUND: $old"
    }
}

$d.Save()
